$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (entry #5) - fill in the previously empty time-log entry
$ws.Range("B11").Value = 43883
$ws.Range("C11").Value = 0.63888888888888895
$ws.Range("D11").Value = 0.68055555555555547
$ws.Range("E11").Value = "-"
$ws.Range("F11").Value = 60
$ws.Range("G11").Value = "Alustasin neljanda  videoga"

# Row 12 (entry #6) - fill in the previously empty time-log entry
$ws.Range("B12").Value = 43884
$ws.Range("C12").Value = 0.76388888888888884
$ws.Range("D12").Value = 0.9375
$ws.Range("E12").Value = 50
$ws.Range("F12").Value = 200
$ws.Range("G12").Value = "Esimesed 6 videot tehtud"
$ws.Range("I12").Value = "x"
$ws.Range("J12").Value = 3

# Update the active selection to match the author's final cursor position
$ws.Range("K12").Select()
